# admin_panel.docx edit
#
# Summary of the change being applied (see the commit's xml diff):
#   1. The stray one-letter "ф" paragraph loses its run (paragraph stays,
#      now empty).
#   2. The "_GoBack" bookmark is removed from inside the
#      "Структура філії [БДБДБДБДБДБДБД]" paragraph...
#   3. ...the final (previously empty) numbered paragraph gets the text
#      "Ф"...
#   4. ...five blank paragraphs are appended...
#   5. ...followed by a paragraph with the "По моему мнению..." text that
#      now carries the relocated "_GoBack" bookmark...
#   6. ...and a last paragraph with the "Передо мной было..." text.
#
# The whole body is rebuilt as one OOXML string and pushed back in with
# Range.InsertXML - the supported way to hand Word a ready-made chunk of
# WordprocessingML instead of poking at runs one at a time.

$d = $word.ActiveDocument

$nsDecls = 'xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se wp14"'

# Paragraphs 1-6 ("Головна" ... "Державні закупівлі") are untouched by the diff.
$headUnchanged = '<w:p w:rsidR="0066663A" w:rsidRPr="00E500B9" w:rsidRDefault="00E500B9" w:rsidP="00E500B9"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Головна</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> [bd]</w:t></w:r></w:p>' + `
    '<w:p w:rsidR="00E500B9" w:rsidRPr="00E500B9" w:rsidRDefault="00E500B9" w:rsidP="00E500B9"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Про ф</w:t></w:r><w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t>ілію</w:t></w:r></w:p>' + `
    '<w:p w:rsidR="00E500B9" w:rsidRPr="00E500B9" w:rsidRDefault="00E500B9" w:rsidP="00E500B9"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t>Історична довідка</w:t></w:r></w:p>' + `
    '<w:p w:rsidR="00E500B9" w:rsidRPr="00E500B9" w:rsidRDefault="00E500B9" w:rsidP="00E500B9"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t>Керівництво</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>[bd]</w:t></w:r></w:p>' + `
    '<w:p w:rsidR="00E500B9" w:rsidRPr="00E500B9" w:rsidRDefault="00E500B9" w:rsidP="00E500B9"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t>Вчена рада</w:t></w:r></w:p>' + `
    '<w:p w:rsidR="00E500B9" w:rsidRPr="00E500B9" w:rsidRDefault="00E500B9" w:rsidP="00E500B9"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t>Державні закупівлі</w:t></w:r></w:p>'

# 1. The "ф" run is dropped - the paragraph (and its pPr) stays, now empty.
$emptyFParagraph = '<w:p w:rsidR="00E500B9" w:rsidRPr="00E500B9" w:rsidRDefault="00E500B9" w:rsidP="00E500B9"><w:pPr><w:rPr><w:lang w:val="uk-UA"/></w:rPr></w:pPr></w:p>'

# 2. Same "Структура філії [БДБДБДБДБДБДБД]" paragraph, minus its bookmark.
$structureParagraph = '<w:p w:rsidR="00E500B9" w:rsidRDefault="00E500B9" w:rsidP="00E500B9"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="uk-UA"/></w:rPr></w:pPr><w:r><w:t>Структура ф</w:t></w:r><w:r w:rsidRPr="00E500B9"><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t>ілії</w:t></w:r><w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>[</w:t></w:r><w:r><w:t>БДБДБДБДБДБДБД</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>]</w:t></w:r></w:p>'

# 3. The final numbered paragraph, now holding the text "Ф".
$filledLastParagraph = '<w:p w:rsidR="00E500B9" w:rsidRPr="00E500B9" w:rsidRDefault="00E500B9" w:rsidP="00E500B9"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="uk-UA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t>Ф</w:t></w:r></w:p>'

# 4. Five new blank paragraphs.
$blankParagraph = '<w:p><w:pPr><w:rPr><w:lang w:val="uk-UA"/></w:rPr></w:pPr></w:p>'
$fiveBlankParagraphs = ""
for ($i = 0; $i -lt 5; $i++) {
    $fiveBlankParagraphs = $fiveBlankParagraphs + $blankParagraph
}

# 5. New paragraph carrying the essay text plus the relocated bookmark.
$gradeText1 = 'По моему мнению, оценка должна ставиться не за сложность какой студент решил то или инное, поставленное перед ним задание, а за путь, котор'
$gradeText2 = 'ым он его решил, пусть он и намного более легкий. '
$gradeParagraph = '<w:p><w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t>' + $gradeText1 + '</w:t></w:r><w:r><w:t xml:space="preserve">' + $gradeText2 + '</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# 6. Final new paragraph.
$taskText = 'Передо мной было поставлено задание сделать все изменяемым «по максимуму», вот я и сделал то, что посчитал целесообразным, так как четкого задания заказчики передо мной не поставили'
$taskParagraph = '<w:p><w:r><w:t>' + $taskText + '</w:t></w:r></w:p>'

$sectPr = '<w:sectPr w:rsidR="00E500B9" w:rsidRPr="00E500B9"><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1134" w:right="850" w:bottom="1134" w:left="1701" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr>'

$body = $headUnchanged + $emptyFParagraph + $structureParagraph + $filledLastParagraph + `
    $fiveBlankParagraphs + $gradeParagraph + $taskParagraph + $sectPr

$documentXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document ' + $nsDecls + '><w:body>' + $body + '</w:body></w:document>'

$packageXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + $documentXml + '</pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($packageXml)
